# edit.ps1 -- apply the "Updated cryptos list" refresh (2024-09-11 GitHub Actions run)
# Rewrites the Price (D) and Volume(1h) (E) columns with the latest scraped
# figures, and swaps the Polygon / VeChain rows (49-50) whose rank order flipped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.644.60"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.320.10"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.34%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "515.75"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "131.92"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("E10").Value = "  -0.10%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.23"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  -1.90%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "23.63"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "2.732.29"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "56.613.89"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "2.326.31"
$ws.Range("E17").Value = "  +0.03%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "10.35"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "328.20"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.46%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.15"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.72"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.67%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "61.07"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -1.17%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "8.61"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.71%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  +1.31%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "167.50"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "0.0₃0718"
$ws.Range("E30").Value = "  -3.75%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.09"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  -0.02%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("E35").Value = "  -1.47%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.94"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.69%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.882"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.98%  "
$ws.Range("E38").Value = "  +0.19%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "38.62"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.73%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "148.63"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("E41").Value = "  -1.45%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.56"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.43%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "276.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  -4.59%  "
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -2.26%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.556"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.52%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "18.28"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0215"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.377"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.52%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "17.06"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
